$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Range("H26").Value = 0.63148
$ws.Range("I26").Value = 0.01859
$ws.Range("H27").Value = 0.11976
$ws.Range("I27").Value = 0.03053
$ws.Range("H28").Value = 0.61068
$ws.Range("I28").Value = 0.01893
$ws.Range("H29").Value = 0.06173
$ws.Range("I29").Value = 0.02748
$ws.Range("H30").Value = 0.63192
$ws.Range("I30").Value = 0.01797
$ws.Range("H31").Value = 0.12076
$ws.Range("I31").Value = 0.03042
$ws.Range("H32").Value = 0.60886
$ws.Range("I32").Value = 0.01752
$ws.Range("H33").Value = 0.06779
$ws.Range("I33").Value = 0.02713
$ws.Range("H34").Value = 0.61627
$ws.Range("I34").Value = 0.01345
$ws.Range("H35").Value = 0.01485
$ws.Range("I35").Value = 0.01437
$ws.Range("H36").Value = 0.61677
$ws.Range("I36").Value = 0.01412
$ws.Range("H37").Value = 0.02126
$ws.Range("I37").Value = 0.01691
$ws.Range("H38").Value = 0.61634
$ws.Range("I38").Value = 0.01344
$ws.Range("H39").Value = 0.01485
$ws.Range("I39").Value = 0.01437
$ws.Range("H40").Value = 0.61667
$ws.Range("I40").Value = 0.01401
$ws.Range("H41").Value = 0.02093
$ws.Range("I41").Value = 0.01708
$ws.Range("H66").Value = 0.57288
$ws.Range("I66").Value = 0.01929
$ws.Range("H67").Value = 0.05028
$ws.Range("I67").Value = 0.02269
$ws.Range("H68").Value = 0.56325
$ws.Range("I68").Value = 0.02078
$ws.Range("H69").Value = 0.0253
$ws.Range("I69").Value = 0.01388
$ws.Range("H70").Value = 0.57436
$ws.Range("I70").Value = 0.01867
$ws.Range("H71").Value = 0.05129
$ws.Range("I71").Value = 0.0216
$ws.Range("H72").Value = 0.56609
$ws.Range("I72").Value = 0.02176
$ws.Range("H73").Value = 0.03136
$ws.Range("I73").Value = 0.01456
$ws.Range("H74").Value = 0.58562
$ws.Range("I74").Value = 0.01371
$ws.Range("H75").Value = 0.01012
$ws.Range("I75").Value = 0.00827
$ws.Range("H76").Value = 0.59026
$ws.Range("I76").Value = 0.01417
$ws.Range("H77").Value = 0.01856
$ws.Range("I77").Value = 0.01147
$ws.Range("H78").Value = 0.58574
$ws.Range("I78").Value = 0.01357
$ws.Range("H79").Value = 0.01012
$ws.Range("I79").Value = 0.00827
$ws.Range("H80").Value = 0.58978
$ws.Range("I80").Value = 0.0141
$ws.Range("H81").Value = 0.01789
$ws.Range("I81").Value = 0.01179
$ws.Range("H106").Value = 0.63469
$ws.Range("I106").Value = 0.02001
$ws.Range("H107").Value = 0.12209
$ws.Range("I107").Value = 0.0321
$ws.Range("H108").Value = 0.61419
$ws.Range("I108").Value = 0.01935
$ws.Range("H109").Value = 0.06542
$ws.Range("I109").Value = 0.02824
$ws.Range("H110").Value = 0.63491
$ws.Range("I110").Value = 0.01941
$ws.Range("H111").Value = 0.12209
$ws.Range("I111").Value = 0.03066
$ws.Range("H112").Value = 0.61076
$ws.Range("I112").Value = 0.02135
$ws.Range("H113").Value = 0.07114
$ws.Range("I113").Value = 0.0281
$ws.Range("H114").Value = 0.6155
$ws.Range("I114").Value = 0.01386
$ws.Range("H115").Value = 0.01788
$ws.Range("I115").Value = 0.00934
$ws.Range("H116").Value = 0.6153
$ws.Range("I116").Value = 0.01385
$ws.Range("H117").Value = 0.02024
$ws.Range("I117").Value = 0.01012
$ws.Range("H118").Value = 0.6155
$ws.Range("I118").Value = 0.01391
$ws.Range("H119").Value = 0.01788
$ws.Range("I119").Value = 0.00934
$ws.Range("H120").Value = 0.61491
$ws.Range("I120").Value = 0.01383
$ws.Range("H121").Value = 0.02024
$ws.Range("I121").Value = 0.01012
$ws.Range("H146").Value = 0.5733
$ws.Range("I146").Value = 0.02183
$ws.Range("H147").Value = 0.04421
$ws.Range("I147").Value = 0.01877
$ws.Range("H148").Value = 0.56488
$ws.Range("I148").Value = 0.02176
$ws.Range("H149").Value = 0.02833
$ws.Range("I149").Value = 0.02061
$ws.Range("H150").Value = 0.57485
$ws.Range("I150").Value = 0.02231
$ws.Range("H151").Value = 0.04522
$ws.Range("I151").Value = 0.01947
$ws.Range("H152").Value = 0.56769
$ws.Range("I152").Value = 0.02609
$ws.Range("H153").Value = 0.04116
$ws.Range("I153").Value = 0.02422
$ws.Range("H154").Value = 0.58495
$ws.Range("I154").Value = 0.01683
$ws.Range("H155").Value = 0.00844
$ws.Range("I155").Value = 0.00925
$ws.Range("H156").Value = 0.58907
$ws.Range("I156").Value = 0.01761
$ws.Range("H157").Value = 0.01824
$ws.Range("I157").Value = 0.01471
$ws.Range("H158").Value = 0.58462
$ws.Range("I158").Value = 0.01708
$ws.Range("H159").Value = 0.00844
$ws.Range("I159").Value = 0.00925
$ws.Range("H160").Value = 0.5884
$ws.Range("I160").Value = 0.01852
$ws.Range("H161").Value = 0.01824
$ws.Range("I161").Value = 0.01471
